# Auto-generated edit script: refresh market-data columns (H-N) across all 8 Leve sheets
# per scheduled market-data runner update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 993.65515
$ws.Range("I33").Value = 808.65216
$ws.Range("K33").Value = 808.65216
$ws.Range("M33").Value = -579.65216
$ws.Range("H51").Value = 41832.668
$ws.Range("J51").Value = 41832.668
$ws.Range("L51").Value = 41832.668
$ws.Range("N51").Value = -42800.668
$ws.Range("H86").Value = 6572.8237
$ws.Range("I86").Value = 6363.857
$ws.Range("K86").Value = 6363.857
$ws.Range("M86").Value = -5240.857
$ws.Range("H89").Value = 6572.8237
$ws.Range("I89").Value = 6363.857
$ws.Range("K89").Value = 31819.285
$ws.Range("M89").Value = -26203.285
$ws.Range("H137").Value = 8876.654
$ws.Range("I137").Value = 8123.3335
$ws.Range("J137").Value = 8974.913
$ws.Range("K137").Value = 24370.0005
$ws.Range("L137").Value = 26924.739
$ws.Range("M137").Value = -21820.0005
$ws.Range("N137").Value = -32024.739
$ws.Range("H141").Value = 7039.7144
$ws.Range("I141").Value = 9822
$ws.Range("K141").Value = 29466
$ws.Range("M141").Value = -24286

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7948735
$ws.Range("I32").Value = 9098319
$ws.Range("K32").Value = 9098319
$ws.Range("M32").Value = -9098032
$ws.Range("H61").Value = 38551588
$ws.Range("I61").Value = 55567576
$ws.Range("K61").Value = 55567576
$ws.Range("M61").Value = -55567364
$ws.Range("H132").Value = 6841.3477
$ws.Range("I132").Value = 3192.0322
$ws.Range("J132").Value = 14383.267
$ws.Range("K132").Value = 9576.096600000001
$ws.Range("L132").Value = 43149.801
$ws.Range("M132").Value = -7046.096600000001
$ws.Range("N132").Value = -48209.801
$ws.Range("H136").Value = 38551588
$ws.Range("I136").Value = 55567576
$ws.Range("K136").Value = 166702728
$ws.Range("M136").Value = -166700178

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1568.9231
$ws.Range("I20").Value = 1213
$ws.Range("J20").Value = 1984.1666
$ws.Range("K20").Value = 1213
$ws.Range("L20").Value = 1984.1666
$ws.Range("M20").Value = -966
$ws.Range("N20").Value = -2478.1666
$ws.Range("H62").Value = 75000
$ws.Range("J62").Value = 75000
$ws.Range("L62").Value = 75000
$ws.Range("N62").Value = -76372
$ws.Range("H65").Value = 75000
$ws.Range("J65").Value = 75000
$ws.Range("L65").Value = 225000
$ws.Range("N65").Value = -231864
$ws.Range("H134").Value = 25681.068
$ws.Range("I134").Value = 3022.1628
$ws.Range("K134").Value = 9066.4884
$ws.Range("M134").Value = -6531.4884
$ws.Range("H105").Value = 3005
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 714415.7
$ws.Range("I31").Value = 11483.625
$ws.Range("K31").Value = 11483.625
$ws.Range("M31").Value = -11188.625
$ws.Range("H34").Value = 714415.7
$ws.Range("I34").Value = 11483.625
$ws.Range("K34").Value = 11483.625
$ws.Range("M34").Value = -11281.625
$ws.Range("H106").Value = 86248
$ws.Range("J106").Value = 86248
$ws.Range("L106").Value = 86248
$ws.Range("N106").Value = -88772
$ws.Range("H122").Value = 2940.2856
$ws.Range("I122").Value = 1807.1111
$ws.Range("K122").Value = 5421.3333
$ws.Range("M122").Value = -2971.3333
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("M65").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 262.5
$ws.Range("J5").Value = 268.5
$ws.Range("L5").Value = 805.5
$ws.Range("N5").Value = -1029.5
$ws.Range("H68").Value = 2119.8823
$ws.Range("J68").Value = 2022.9333
$ws.Range("L68").Value = 6068.7999
$ws.Range("N68").Value = -7690.7999
$ws.Range("H71").Value = 2119.8823
$ws.Range("J71").Value = 2022.9333
$ws.Range("L71").Value = 18206.3997
$ws.Range("N71").Value = -26318.3997
$ws.Range("H121").Value = 2768.182
$ws.Range("I121").Value = 1374
$ws.Range("J121").Value = 5208
$ws.Range("K121").Value = 4122
$ws.Range("L121").Value = 15624
$ws.Range("M121").Value = -2812
$ws.Range("N121").Value = -18244
$ws.Range("H131").Value = 6048.1816
$ws.Range("J131").Value = 1998.3334
$ws.Range("L131").Value = 5995.0002
$ws.Range("N131").Value = -16075.0002
$ws.Range("H135").Value = 262.5
$ws.Range("J135").Value = 268.5
$ws.Range("L135").Value = 2416.5
$ws.Range("N135").Value = -7486.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2872.5715
$ws.Range("I80").Value = 2651.25
$ws.Range("J80").Value = 3167.6667
$ws.Range("K80").Value = 2651.25
$ws.Range("L80").Value = 3167.6667
$ws.Range("M80").Value = -1653.25
$ws.Range("N80").Value = -5163.6667
$ws.Range("H83").Value = 2872.5715
$ws.Range("I83").Value = 2651.25
$ws.Range("J83").Value = 3167.6667
$ws.Range("K83").Value = 13256.25
$ws.Range("L83").Value = 15838.3335
$ws.Range("M83").Value = -8264.25
$ws.Range("N83").Value = -25822.3335
$ws.Range("H104").Value = 100000
$ws.Range("J104").Value = 100000
$ws.Range("L104").Value = 100000
$ws.Range("N104").Value = -106988
$ws.Range("H126").Value = 3781.2666
$ws.Range("I126").Value = 2845.5715
$ws.Range("K126").Value = 8536.7145
$ws.Range("M126").Value = -6066.7145
$ws.Range("H132").Value = 47622680
$ws.Range("I132").Value = 55559330
$ws.Range("K132").Value = 166677990
$ws.Range("M132").Value = -166675460
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 115099.78
$ws.Range("I7").Value = 3577.8
$ws.Range("J7").Value = 254502.25
$ws.Range("K7").Value = 3577.8
$ws.Range("L7").Value = 254502.25
$ws.Range("M7").Value = -3465.8
$ws.Range("N7").Value = -254726.25
$ws.Range("H22").Value = 3258.5454
$ws.Range("I22").Value = 4220.2856
$ws.Range("J22").Value = 1575.5
$ws.Range("K22").Value = 4220.2856
$ws.Range("L22").Value = 1575.5
$ws.Range("M22").Value = -3925.2856
$ws.Range("N22").Value = -2165.5
$ws.Range("H27").Value = 3258.5454
$ws.Range("I27").Value = 4220.2856
$ws.Range("J27").Value = 1575.5
$ws.Range("K27").Value = 4220.2856
$ws.Range("L27").Value = 1575.5
$ws.Range("M27").Value = -4113.2856
$ws.Range("N27").Value = -1789.5
$ws.Range("H68").Value = 3199.75
$ws.Range("I68").Value = 2499.5
$ws.Range("J68").Value = 3900
$ws.Range("K68").Value = 2499.5
$ws.Range("L68").Value = 3900
$ws.Range("M68").Value = -1750.5
$ws.Range("N68").Value = -5398
$ws.Range("H71").Value = 3199.75
$ws.Range("I71").Value = 2499.5
$ws.Range("J71").Value = 3900
$ws.Range("K71").Value = 12497.5
$ws.Range("L71").Value = 19500
$ws.Range("M71").Value = -8753.5
$ws.Range("N71").Value = -26988
$ws.Range("H100").Value = 6022.923
$ws.Range("I100").Value = 2659.8
$ws.Range("K100").Value = 2659.8
$ws.Range("M100").Value = -2118.8
$ws.Range("H126").Value = 115099.78
$ws.Range("I126").Value = 3577.8
$ws.Range("J126").Value = 254502.25
$ws.Range("K126").Value = 10733.4
$ws.Range("L126").Value = 763506.75
$ws.Range("M126").Value = -8263.400000000001
$ws.Range("N126").Value = -768446.75
$ws.Range("H61").Value = 2586
$ws.Range("I61").Value = 2586
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2586
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2384
$ws.Range("H113").Value = 2586
$ws.Range("I113").Value = 2586
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2586
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -416
$ws.Range("N61").ClearContents()
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 57983.715
$ws.Range("I81").Value = 40577.4
$ws.Range("K81").Value = 81154.8
$ws.Range("M81").Value = -80093.8
$ws.Range("H84").Value = 57983.715
$ws.Range("I84").Value = 40577.4
$ws.Range("K84").Value = 405774
$ws.Range("M84").Value = -400470
$ws.Range("H107").Value = 71429976
$ws.Range("I107").Value = 100001780
$ws.Range("K107").Value = 300005340
$ws.Range("M107").Value = -300003420
$ws.Range("H126").Value = 3325.4897
$ws.Range("J126").Value = 4327.4375
$ws.Range("L126").Value = 12982.3125
$ws.Range("N126").Value = -17922.3125
$ws.Range("H132").Value = 274494.9
$ws.Range("I132").Value = 2750.1924
$ws.Range("K132").Value = 8250.5772
$ws.Range("M132").Value = -5720.5772
$ws.Range("H92").Value = 55000
$ws.Range("J92").Value = 55000
$ws.Range("L92").Value = 55000
$ws.Range("N92").Value = -59992
